$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "abhi"
$ws.Range("C3").Value = "08:06:29"

# "2025-01-16" looks like a date, so Excel would normally auto-convert it to
# a date serial number. Force it to be stored as literal text instead, then
# strip the resulting formatting so the cell ends up unstyled (like the rest
# of the data row).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-01-16"
$ws.Range("B3").ClearFormats()
